# Auto-generated script applying the cryptos.xlsx price/percentage/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.870.40"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.429.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "411.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.26%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.733"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "44.52"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000222"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +8.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.44"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +8.17%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.971.00"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.141"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.47"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.411.33"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.55"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +7.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.10"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "61.967.73"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "499.81"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +50.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "92.86"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.38"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.35"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.32"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "35.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +15.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.12"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +10.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.78"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.72"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.78"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.24"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.16%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.115"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.169"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.81"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.48%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0506"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.91"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.48"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.34%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.136"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.78%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.69"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +15.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.319"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.40"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +11.05%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "145.46"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.07"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.22%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.25"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +7.03%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "120.67"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +41.28%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.143"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +15.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.13"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +7.60%  "
